$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "https://www.quora.com/What-is-the-software-architecture-of-Twitter"
$ws.Range("D2").Value = "not the length"
$ws.Range("D7").Value = "nope"

$ws.Range("C7").Select()
